$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 0.0027282
$ws.Range("E2").Value = 0.3725572890751533
$ws.Range("G2").Value = 0.2494892361374987
$ws.Range("I2").Value = 0.3669021
$ws.Range("L2").Value = 0.5961429402307628
$ws.Range("M2").Value = 0.08239116666666667
$ws.Range("N2").Value = 12.90833032859821
$ws.Range("O2").Value = 3.475337169561454

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 0.06802560000000066
$ws.Range("E2").Value = 0.3710309236677573
$ws.Range("I2").Value = 0.7381097563895144
$ws.Range("L2").Value = 0.3168957764359041
$ws.Range("M2").Value = 0.08173224999999999
$ws.Range("N2").Value = 8.975002603116078
$ws.Range("O2").Value = 4.340088225780883

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 0.09762676620383053
$ws.Range("B2").Value = 0.02931738907515262
$ws.Range("E2").Value = 0.1495921412480008
$ws.Range("I2").Value = 0.2469301557893744
$ws.Range("M2").Value = 0.04645608333333337
$ws.Range("N2").Value = 9.012664873503194
$ws.Range("O2").Value = 4.788897617609521
